$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I20").Value = -0.3647414710054033
$ws.Range("J20").Value = 0.2426200462448084
$ws.Range("K20").Value = 0.2814819410217358
$ws.Range("L20").Value = 2.570540023702133
